$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ReadOldPrices" step: shift the previous "Prices" (col B) values into
# "Old Prices" (col D), then refresh "Prices" (col B) and "Euro" (col C)
# with newly fetched values.

# Row 2 (Binance)
$ws.Range("D2").Value2 = $ws.Range("B2").Value2
$ws.Range("B2").Value2 = 46217.93
$ws.Range("C2").Value2 = 40692.22

# Row 3 (KuCoin)
$ws.Range("D3").Value2 = $ws.Range("B3").Value2
$ws.Range("B3").Value2 = 46215.1
$ws.Range("C3").Value2 = 40689.73

# Row 4 (Coinbase)
$ws.Range("D4").Value2 = $ws.Range("B4").Value2
$ws.Range("B4").Value2 = 46192.2
$ws.Range("C4").Value2 = 40669.56

# Row 5 (CMC)
$ws.Range("D5").Value2 = $ws.Range("B5").Value2
$ws.Range("B5").Value2 = 46197.55
$ws.Range("C5").Value2 = 40674.28

# Row 6 (CoinGecko) - B/D hold text (string) price values, not numbers.
# Use a leading apostrophe to force text entry (preserves the trailing
# space in "46276.50 "/"46131.60 " instead of Excel auto-converting it to
# a number), then reset the cell style so no new number format sticks.
$ws.Range("D6").Formula = "'" + $ws.Range("B6").Value2
$ws.Range("D6").Style = "Normal"

$ws.Range("B6").Formula = "'46276.50" + [char]0x00A0
$ws.Range("B6").Style = "Normal"

$ws.Range("C6").Value2 = 40743.79

$wb.Save()
